$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2098445595854922
$ws.Range("C2").Value = 0.5440414507772021
$ws.Range("J2").Value = 0.0155440414507772
$ws.Range("P2").Value = 0.1450777202072539
$ws.Range("S2").Value = 0.08549222797927461

$ws.Range("B3").Value = 0.00909090909090909
$ws.Range("C3").Value = 0.03636363636363636
$ws.Range("J3").Value = 0.04090909090909091
$ws.Range("P3").Value = 0.7
$ws.Range("S3").Value = 0.2136363636363636

$ws.Range("P4").Value = 0.675
$ws.Range("S4").Value = 0.325

$ws.Range("B6").Value = 0.07547169811320754
$ws.Range("D6").Value = 0.01415094339622642
$ws.Range("F6").Value = 0.04716981132075472
$ws.Range("J6").Value = 0.2358490566037736
$ws.Range("O6").Value = 0.01886792452830189
$ws.Range("Q6").Value = 0.1415094339622641
$ws.Range("R6").Value = 0.09433962264150944
$ws.Range("S6").Value = 0.3726415094339622

$ws.Range("B7").Value = 0.1022727272727273
$ws.Range("D7").Value = 0.005681818181818182
$ws.Range("F7").Value = 0.03409090909090909
$ws.Range("J7").Value = 0.2102272727272727
$ws.Range("O7").Value = 0.005681818181818182
$ws.Range("Q7").Value = 0.1875
$ws.Range("R7").Value = 0.06818181818181818
$ws.Range("S7").Value = 0.3863636363636364

$ws.Range("B8").Value = 0.1348600508905853
$ws.Range("D8").Value = 0.02544529262086514
$ws.Range("F8").Value = 0.09923664122137404
$ws.Range("J8").Value = 0.1195928753180662
$ws.Range("O8").Value = 0.01272264631043257
$ws.Range("Q8").Value = 0.1628498727735369
$ws.Range("R8").Value = 0.08651399491094147
$ws.Range("S8").Value = 0.3587786259541985

$ws.Range("B9").Value = 0.06914893617021277
$ws.Range("D9").Value = 0.005319148936170213
$ws.Range("E9").Value = 0.005319148936170213
$ws.Range("F9").Value = 0.05851063829787234
$ws.Range("J9").Value = 0.1436170212765958
$ws.Range("O9").Value = 0.02127659574468085
$ws.Range("Q9").Value = 0.2021276595744681
$ws.Range("R9").Value = 0.07446808510638298
$ws.Range("S9").Value = 0.4202127659574468

$ws.Range("B10").Value = 0.1427546628407461
$ws.Range("D10").Value = 0.01865136298421808
$ws.Range("E10").Value = 0.002869440459110474
$ws.Range("F10").Value = 0.06384505021520803
$ws.Range("J10").Value = 0.1262553802008608
$ws.Range("O10").Value = 0.0157819225251076
$ws.Range("Q10").Value = 0.2087517934002869
$ws.Range("R10").Value = 0.08249641319942611
$ws.Range("S10").Value = 0.3385939741750359

$ws.Range("G11").Value = 0.1654929577464789
$ws.Range("J11").Value = 0.08450704225352113
$ws.Range("K11").Value = 0.2007042253521127
$ws.Range("L11").Value = 0.5422535211267606
$ws.Range("S11").Value = 0.007042253521126761

$ws.Range("G12").Value = 0.6875
$ws.Range("J12").Value = 0.26875
$ws.Range("L12").Value = 0.01875
$ws.Range("S12").Value = 0.025

$ws.Range("F13").Value = 0.02857142857142857
$ws.Range("G13").Value = 0.6857142857142857
$ws.Range("J13").Value = 0.2285714285714286
$ws.Range("S13").Value = 0.05714285714285714

$ws.Range("F15").Value = 0.01507537688442211
$ws.Range("H15").Value = 0.1306532663316583
$ws.Range("I15").Value = 0.08040201005025126
$ws.Range("J15").Value = 0.4020100502512563
$ws.Range("K15").Value = 0.04522613065326633
$ws.Range("M15").Value = 0.02010050251256281
$ws.Range("O15").Value = 0.03517587939698492
$ws.Range("S15").Value = 0.271356783919598

$ws.Range("F16").Value = 0.02164502164502164
$ws.Range("H16").Value = 0.1298701298701299
$ws.Range("I16").Value = 0.08658008658008658
$ws.Range("J16").Value = 0.4415584415584415
$ws.Range("K16").Value = 0.09523809523809523
$ws.Range("M16").Value = 0.02164502164502164
$ws.Range("O16").Value = 0.02164502164502164
$ws.Range("S16").Value = 0.1818181818181818

$ws.Range("H17").Value = 0.168141592920354
$ws.Range("I17").Value = 0.07743362831858407
$ws.Range("J17").Value = 0.4756637168141593
$ws.Range("K17").Value = 0.08849557522123894
$ws.Range("M17").Value = 0.01106194690265487
$ws.Range("O17").Value = 0.04424778761061947
$ws.Range("S17").Value = 0.1349557522123894

$ws.Range("F18").Value = 0.01036269430051814
$ws.Range("H18").Value = 0.1865284974093264
$ws.Range("I18").Value = 0.1036269430051813
$ws.Range("J18").Value = 0.3989637305699482
$ws.Range("K18").Value = 0.07772020725388601
$ws.Range("M18").Value = 0.0155440414507772
$ws.Range("O18").Value = 0.08808290155440414
$ws.Range("S18").Value = 0.1191709844559585

$ws.Range("F19").Value = 0.01449275362318841
$ws.Range("H19").Value = 0.1819645732689211
$ws.Range("I19").Value = 0.07971014492753623
$ws.Range("J19").Value = 0.4082125603864734
$ws.Range("K19").Value = 0.1111111111111111
$ws.Range("M19").Value = 0.01610305958132045
$ws.Range("N19").Value = 0.0008051529790660225
$ws.Range("O19").Value = 0.07165861513687601
$ws.Range("S19").Value = 0.1159420289855072

Write-Output "applied team specific time data changes"